# Updated cryptos list refresh: write new Price (column D) and
# Volume(1h) (column E) values for the rows that changed.
#
# Price cells that look like plain numbers (e.g. "211.29") are forced
# to remain plain text by temporarily switching the cell to Text
# number format before assigning the value, then restoring the
# default "Normal" style so no stray formatting is left behind.
# Price values that already contain extra separators (e.g.
# "27.856.80") or non-numeric characters are never auto-converted by
# Excel, so they can be assigned directly. Volume(1h) strings contain
# surrounding spaces/percent signs and are likewise never numeric, so
# they are always assigned directly as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.856.80'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0881'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.861.02'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.621.31'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '27.865.72'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = '0.0₃0718'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.65%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').Value = '1.413.60'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.996'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0170'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.554'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.854'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  -2.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '1.770.27'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('E46').Value = '  -3.95%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0503'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.995'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
